$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: single-cell corrections to column C (covid_deaths counts) ---
$ws.Cells.Item(33, 3).Value = 3
$ws.Cells.Item(65, 3).Value = 6
$ws.Cells.Item(1005, 3).Value = 46
$ws.Cells.Item(1042, 3).Value = 48
$ws.Cells.Item(1082, 3).Value = 19
$ws.Cells.Item(1118, 3).Value = 25
$ws.Cells.Item(1131, 3).Value = 70
$ws.Cells.Item(1187, 3).Value = 31
$ws.Cells.Item(1203, 3).Value = 11
$ws.Cells.Item(1210, 3).Value = 38
$ws.Cells.Item(1226, 3).Value = 2
$ws.Cells.Item(1236, 3).Value = 38
$ws.Cells.Item(1266, 3).Value = 31
$ws.Cells.Item(1284, 3).Value = 43
$ws.Cells.Item(1325, 3).Value = 14
$ws.Cells.Item(1337, 3).Value = 8
$ws.Cells.Item(1339, 3).Value = 27
$ws.Cells.Item(1344, 3).Value = 23
$ws.Cells.Item(1377, 3).Value = 21
$ws.Cells.Item(1394, 3).Value = 16
$ws.Cells.Item(1398, 3).Value = 13
$ws.Cells.Item(1399, 3).Value = 5
$ws.Cells.Item(1406, 3).Value = 13

# --- Part 2: rows 1409-1424 shift (re-aggregated tail of the series) and new rows 1425-1434 appended ---
$ws.Cells.Item(1409, 2).Value = "50-59"
$ws.Cells.Item(1409, 3).Value = 1
$ws.Cells.Item(1410, 2).Value = "60-69"
$ws.Cells.Item(1410, 3).Value = 3
$ws.Cells.Item(1411, 2).Value = "70-79"
$ws.Cells.Item(1411, 3).Value = 3
$ws.Cells.Item(1412, 1).Value = 44240
$ws.Cells.Item(1412, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1412, 2).Value = "80+"
$ws.Cells.Item(1412, 3).Value = 7
$ws.Cells.Item(1413, 2).Value = "40-49"
$ws.Cells.Item(1413, 3).Value = 2
$ws.Cells.Item(1414, 2).Value = "50-59"
$ws.Cells.Item(1414, 3).Value = 1
$ws.Cells.Item(1415, 2).Value = "60-69"
$ws.Cells.Item(1415, 3).Value = 3
$ws.Cells.Item(1416, 2).Value = "70-79"
$ws.Cells.Item(1416, 3).Value = 2
$ws.Cells.Item(1417, 1).Value = 44241
$ws.Cells.Item(1417, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1417, 2).Value = "80+"
$ws.Cells.Item(1417, 3).Value = 4
$ws.Cells.Item(1418, 2).Value = "50-59"
$ws.Cells.Item(1418, 3).Value = 1
$ws.Cells.Item(1419, 2).Value = "60-69"
$ws.Cells.Item(1419, 3).Value = 5
$ws.Cells.Item(1420, 2).Value = "70-79"
$ws.Cells.Item(1420, 3).Value = 5
$ws.Cells.Item(1421, 1).Value = 44242
$ws.Cells.Item(1421, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1421, 2).Value = "80+"
$ws.Cells.Item(1421, 3).Value = 12
$ws.Cells.Item(1422, 2).Value = "50-59"
$ws.Cells.Item(1422, 3).Value = 2
$ws.Cells.Item(1423, 2).Value = "60-69"
$ws.Cells.Item(1423, 3).Value = 4
$ws.Cells.Item(1424, 2).Value = "70-79"
$ws.Cells.Item(1424, 3).Value = 3
$ws.Cells.Item(1425, 1).Value = 44243
$ws.Cells.Item(1425, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1425, 2).Value = "80+"
$ws.Cells.Item(1425, 3).Value = 5
$ws.Cells.Item(1426, 1).Value = 44244
$ws.Cells.Item(1426, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1426, 2).Value = "40-49"
$ws.Cells.Item(1426, 3).Value = 2
$ws.Cells.Item(1427, 1).Value = 44244
$ws.Cells.Item(1427, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1427, 2).Value = "60-69"
$ws.Cells.Item(1427, 3).Value = 3
$ws.Cells.Item(1428, 1).Value = 44244
$ws.Cells.Item(1428, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1428, 2).Value = "70-79"
$ws.Cells.Item(1428, 3).Value = 3
$ws.Cells.Item(1429, 1).Value = 44244
$ws.Cells.Item(1429, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1429, 2).Value = "80+"
$ws.Cells.Item(1429, 3).Value = 3
$ws.Cells.Item(1430, 1).Value = 44245
$ws.Cells.Item(1430, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1430, 2).Value = "60-69"
$ws.Cells.Item(1430, 3).Value = 2
$ws.Cells.Item(1431, 1).Value = 44245
$ws.Cells.Item(1431, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1431, 2).Value = "70-79"
$ws.Cells.Item(1431, 3).Value = 4
$ws.Cells.Item(1432, 1).Value = 44245
$ws.Cells.Item(1432, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1432, 2).Value = "80+"
$ws.Cells.Item(1432, 3).Value = 1
$ws.Cells.Item(1433, 1).Value = 44246
$ws.Cells.Item(1433, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1433, 2).Value = "60-69"
$ws.Cells.Item(1433, 3).Value = 1
$ws.Cells.Item(1434, 1).Value = 44246
$ws.Cells.Item(1434, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1434, 2).Value = "80+"
$ws.Cells.Item(1434, 3).Value = 1
